# tools.ps1 now has Multi-Config support, Setup Guide Updated based on
# feedback, Connection_HID should be complete now.
#
# Removes the now-obsolete "Test whether the current configurations are
# standalone ..." bullet (and its "Do not require vcpkg ..." sub-bullet)
# plus the trailing empty bullet paragraph that followed them, since the
# Connection_HID setup guide work they referenced is complete.

$d = $word.ActiveDocument

# Locate the last paragraph whose text contains the anchor phrase - this
# is the bullet that should remain as the new end of the list/document.
$anchorText = "This can be done because ps1 script files can be converted into exe files"
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*$anchorText*") {
        $anchorIndex = $i
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not find anchor paragraph: $anchorText"
}

# Everything from the paragraph right after the anchor through the end of
# the document (the two stray bullets plus the trailing empty paragraph)
# gets removed.
if ($anchorIndex -lt $d.Paragraphs.Count) {
    $startPara = $d.Paragraphs($anchorIndex + 1)
    $lastPara = $d.Paragraphs($d.Paragraphs.Count)
    $delRange = $d.Range($startPara.Range.Start, $lastPara.Range.End)
    $delRange.Delete()
}
